$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 320.375
$ws.Range("J33").Value = 455.7
$ws.Range("L33").Value = 455.7
$ws.Range("N33").Value = -913.7
$ws.Range("H43").Value = 3349.5
$ws.Range("I43").Value = 2666.3333
$ws.Range("K43").Value = 2666.3333
$ws.Range("M43").Value = -2597.3333
$ws.Range("H55").Value = 125
$ws.Range("I55").Value = 137.5
$ws.Range("K55").Value = 137.5
$ws.Range("M55").Value = 76.5
$ws.Range("H115").Value = 900
$ws.Range("I115").Value = 900
$ws.Range("K115").Value = 2700
$ws.Range("M115").Value = -1133
$ws.Range("H138").Value = 2005.5217
$ws.Range("I138").Value = 1390.6364
$ws.Range("J138").Value = 3566.3845
$ws.Range("K138").Value = 4171.9092
$ws.Range("L138").Value = 10699.1535
$ws.Range("M138").Value = 968.0907999999999
$ws.Range("N138").Value = -20979.1535
$ws.Range("H141").Value = 1991.091
$ws.Range("I141").Value = 1991.091
$ws.Range("K141").Value = 5973.272999999999
$ws.Range("M141").Value = -793.2729999999992

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16000.448
$ws.Range("I32").Value = 16227.728
$ws.Range("K32").Value = 16227.728
$ws.Range("M32").Value = -15940.728
$ws.Range("H61").Value = 14690.223
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 16151.5
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 16151.5
$ws.Range("M61").Value = -2788
$ws.Range("N61").Value = -16575.5
$ws.Range("H74").Value = 382516
$ws.Range("I74").Value = 462600.94
$ws.Range("J74").Value = 35481.332
$ws.Range("K74").Value = 462600.94
$ws.Range("L74").Value = 35481.332
$ws.Range("M74").Value = -461726.94
$ws.Range("N74").Value = -37229.332
$ws.Range("H77").Value = 382516
$ws.Range("I77").Value = 462600.94
$ws.Range("J77").Value = 35481.332
$ws.Range("K77").Value = 2313004.7
$ws.Range("L77").Value = 177406.66
$ws.Range("M77").Value = -2308636.7
$ws.Range("N77").Value = -186142.66
$ws.Range("H102").Value = 3207.4443
$ws.Range("I102").Value = 3136.5715
$ws.Range("K102").Value = 3136.5715
$ws.Range("M102").Value = -1514.5715
$ws.Range("H122").Value = 2295.75
$ws.Range("I122").Value = 1901.4615
$ws.Range("K122").Value = 5704.3845
$ws.Range("M122").Value = -3254.3845
$ws.Range("H132").Value = 1414.7894
$ws.Range("I132").Value = 992.5625
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 2977.6875
$ws.Range("L132").Value = 11000.0001
$ws.Range("M132").Value = -447.6875
$ws.Range("N132").Value = -16060.0001
$ws.Range("H136").Value = 14690.223
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 16151.5
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 48454.5
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -53554.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 590
$ws.Range("I22").Value = 590
$ws.Range("K22").Value = 590
$ws.Range("M22").Value = -417
$ws.Range("H94").Value = 5084.885
$ws.Range("I94").Value = 5218.5454
$ws.Range("K94").Value = 5218.5454
$ws.Range("M94").Value = -4767.5454
$ws.Range("H105").Value = 4395.222
$ws.Range("I105").Value = 4302.6
$ws.Range("J105").Value = 4511
$ws.Range("K105").Value = 4302.6
$ws.Range("L105").Value = 4511
$ws.Range("M105").Value = -2555.6
$ws.Range("N105").Value = -8005
$ws.Range("H107").Value = 4313.125
$ws.Range("I107").Value = 4516.1665
$ws.Range("J107").Value = 3704
$ws.Range("K107").Value = 4516.1665
$ws.Range("L107").Value = 3704
$ws.Range("M107").Value = -2596.1665
$ws.Range("N107").Value = -7544

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 141.72728
$ws.Range("I7").Value = 72.5
$ws.Range("J7").Value = 224.8
$ws.Range("K7").Value = 72.5
$ws.Range("L7").Value = 224.8
$ws.Range("M7").Value = 40.5
$ws.Range("N7").Value = -450.8
$ws.Range("H31").Value = 3704481.2
$ws.Range("I31").Value = 3704481.2
$ws.Range("K31").Value = 3704481.2
$ws.Range("M31").Value = -3704186.2
$ws.Range("H34").Value = 3704481.2
$ws.Range("I34").Value = 3704481.2
$ws.Range("K34").Value = 3704481.2
$ws.Range("M34").Value = -3704279.2
$ws.Range("H50").Value = 36315.832
$ws.Range("J50").Value = 36315.832
$ws.Range("L50").Value = 36315.832
$ws.Range("N50").Value = -37565.832
$ws.Range("H58").Value = 1307.4231
$ws.Range("I58").Value = 1128.0526
$ws.Range("K58").Value = 1128.0526
$ws.Range("M58").Value = -925.0526
$ws.Range("H99").Value = 8536.5
$ws.Range("J99").Value = 7601.75
$ws.Range("L99").Value = 7601.75
$ws.Range("N99").Value = -10597.75
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H126").Value = 8536.5
$ws.Range("J126").Value = 7601.75
$ws.Range("L126").Value = 22805.25
$ws.Range("N126").Value = -27745.25
$ws.Range("H134").Value = 2860.724
$ws.Range("I134").Value = 2294
$ws.Range("J134").Value = 5033.1665
$ws.Range("K134").Value = 6882
$ws.Range("L134").Value = 15099.4995
$ws.Range("M134").Value = -4347
$ws.Range("N134").Value = -20169.4995
$ws.Range("H136").Value = 1307.4231
$ws.Range("I136").Value = 1128.0526
$ws.Range("K136").Value = 3384.1578
$ws.Range("M136").Value = -834.1578

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H107").Value = 1816.9286
$ws.Range("I107").Value = 2543.375
$ws.Range("J107").Value = 848.3333
$ws.Range("K107").Value = 7630.125
$ws.Range("L107").Value = 2544.9999
$ws.Range("M107").Value = -5710.125
$ws.Range("N107").Value = -6384.9999
$ws.Range("H132").Value = 1838.5
$ws.Range("J132").Value = 1931.6666
$ws.Range("L132").Value = 17384.9994
$ws.Range("N132").Value = -22444.9994
$ws.Range("H140").Value = 1921.875
$ws.Range("I140").Value = 1921.875
$ws.Range("K140").Value = 5765.625
$ws.Range("M140").Value = -585.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 139.2
$ws.Range("I2").Value = 91.75
$ws.Range("J2").Value = 234.1
$ws.Range("K2").Value = 91.75
$ws.Range("L2").Value = 234.1
$ws.Range("M2").Value = 21.25
$ws.Range("N2").Value = -460.1
$ws.Range("H97").Value = 1243.875
$ws.Range("I97").Value = 1113.5
$ws.Range("J97").Value = 1635
$ws.Range("K97").Value = 1113.5
$ws.Range("L97").Value = 1635
$ws.Range("M97").Value = -617.5
$ws.Range("N97").Value = -2627
$ws.Range("H113").Value = 2830.8635
$ws.Range("I113").Value = 2572.3333
$ws.Range("J113").Value = 3384.8572
$ws.Range("K113").Value = 2572.3333
$ws.Range("L113").Value = 3384.8572
$ws.Range("M113").Value = -402.3332999999998
$ws.Range("N113").Value = -7724.8572
$ws.Range("H122").Value = 4098.75
$ws.Range("I122").Value = 4209.4546
$ws.Range("J122").Value = 3963.4443
$ws.Range("K122").Value = 12628.3638
$ws.Range("L122").Value = 11890.3329
$ws.Range("M122").Value = -10178.3638
$ws.Range("N122").Value = -16790.3329
$ws.Range("H126").Value = 2759.6
$ws.Range("I126").Value = 1949.5
$ws.Range("K126").Value = 5848.5
$ws.Range("M126").Value = -3378.5
$ws.Range("H139").Value = 104081.375
$ws.Range("J139").Value = 104081.375
$ws.Range("L139").Value = 104081.375
$ws.Range("N139").Value = -114361.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1912.7354
$ws.Range("I55").Value = 1311.2727
$ws.Range("K55").Value = 1311.2727
$ws.Range("M55").Value = -1138.2727
$ws.Range("H68").Value = 4050.4614
$ws.Range("I68").Value = 3379.8572
$ws.Range("K68").Value = 3379.8572
$ws.Range("M68").Value = -2630.8572
$ws.Range("H71").Value = 4050.4614
$ws.Range("I71").Value = 3379.8572
$ws.Range("K71").Value = 16899.286
$ws.Range("M71").Value = -13155.286
$ws.Range("H93").Value = 1271.5454
$ws.Range("I93").Value = 1271.5454
$ws.Range("K93").Value = 1271.5454
$ws.Range("M93").Value = -23.54539999999997

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4003
$ws.Range("J62").Value = 4003
$ws.Range("L62").Value = 4003
$ws.Range("N62").Value = -5251
$ws.Range("H65").Value = 4003
$ws.Range("J65").Value = 4003
$ws.Range("L65").Value = 20015
$ws.Range("N65").Value = -26255
$ws.Range("H74").Value = 20749.5
$ws.Range("I74").Value = 20000
$ws.Range("K74").Value = 20000
$ws.Range("M74").Value = -19064
$ws.Range("H77").Value = 20749.5
$ws.Range("I77").Value = 20000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55320
$ws.Range("H93").Value = 100000
$ws.Range("J93").Value = 100000
$ws.Range("L93").Value = 100000
$ws.Range("N93").Value = -104992
$ws.Range("H136").Value = 15329.208
$ws.Range("I136").Value = 19404.445
$ws.Range("J136").Value = 3103.5
$ws.Range("K136").Value = 58213.335
$ws.Range("L136").Value = 9310.5
$ws.Range("M136").Value = -14410.5
